$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new value is a "clean" number (e.g. "0.477") need
# to be pinned to Text format first, otherwise Excel auto-converts the
# assigned string into a numeric value (losing the fixed-precision text,
# e.g. "0.477" -> 0.476999999999...). Two-dot values like "25.868.65" or
# subscript-containing values like "0.0₃0710" already fail numeric parsing
# so they naturally stay text and do not need this treatment.
$textForcedCells = @(
    "D5", "D7", "D9", "D10", "D11", "D17", "D20", "D22", "D25", "D27", "D29", "D38", "D39", "D40", "D41", "D43", "D44", "D46", "D48", "D51"
)
foreach ($cellRef in $textForcedCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "25.868.65"
$ws.Range("E2").Value = "  -0.43%  "
$ws.Range("D3").Value = "1.598.78"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "208.81"
$ws.Range("E5").Value = "  -2.38%  "
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").Value = "0.477"
$ws.Range("E7").Value = "  -5.34%  "
$ws.Range("E8").Value = "  -2.59%  "
$ws.Range("D9").Value = "0.0611"
$ws.Range("D10").Value = "17.76"
$ws.Range("E10").Value = "  -3.95%  "
$ws.Range("D11").Value = "0.0785"
$ws.Range("E11").Value = "  -0.84%  "
$ws.Range("D12").Value = "1.822.06"
$ws.Range("E12").Value = "  -1.99%  "
$ws.Range("D13").Value = "1.602.60"
$ws.Range("E13").Value = "  -2.24%  "
$ws.Range("E14").Value = "  -3.63%  "
$ws.Range("E15").Value = "  -4.44%  "
$ws.Range("D16").Value = "25.861.85"
$ws.Range("E16").Value = "  -0.46%  "
$ws.Range("D17").Value = "60.38"
$ws.Range("E17").Value = "  -2.09%  "
$ws.Range("D18").Value = "0.0₃0710"
$ws.Range("E18").Value = "  -4.73%  "
$ws.Range("E19").Value = "  +0.04%  "
$ws.Range("D20").Value = "189.00"
$ws.Range("E20").Value = "  -0.63%  "
$ws.Range("E21").Value = "  -1.78%  "
$ws.Range("D22").Value = "9.30"
$ws.Range("E22").Value = "  -2.82%  "
$ws.Range("E23").Value = "  -3.25%  "
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("D25").Value = "141.61"
$ws.Range("E25").Value = "  -1.31%  "
$ws.Range("E26").Value = "  -3.60%  "
$ws.Range("D27").Value = "1.69"
$ws.Range("E27").Value = "  -3.71%  "
$ws.Range("E28").Value = "  -4.07%  "
$ws.Range("D29").Value = "14.89"
$ws.Range("E29").Value = "  -1.99%  "
$ws.Range("E30").Value = "  -2.43%  "
$ws.Range("E31").Value = "  -4.57%  "
$ws.Range("E32").Value = "  -2.71%  "
$ws.Range("E33").Value = "  -5.00%  "
$ws.Range("E34").Value = "  -1.00%  "
$ws.Range("E35").Value = "  -1.93%  "
$ws.Range("D36").Value = "1.103.76"
$ws.Range("E36").Value = "  -2.45%  "
$ws.Range("E37").Value = "  -3.00%  "
$ws.Range("D38").Value = "0.795"
$ws.Range("E38").Value = "  -8.14%  "
$ws.Range("D39").Value = "0.0150"
$ws.Range("E39").Value = "  -2.99%  "
$ws.Range("D40").Value = "0.495"
$ws.Range("E40").Value = "  -5.55%  "
$ws.Range("D41").Value = "95.48"
$ws.Range("E41").Value = "  -3.05%  "
$ws.Range("D42").Value = "1.734.47"
$ws.Range("E42").Value = "  -1.97%  "
$ws.Range("D43").Value = "5.05"
$ws.Range("E43").Value = "  -3.62%  "
$ws.Range("D44").Value = "0.739"
$ws.Range("E44").Value = "  -5.05%  "
$ws.Range("D45").Value = "0.0⁦0104"
$ws.Range("E45").Value = "  -8.55%  "
$ws.Range("D46").Value = "52.94"
$ws.Range("E46").Value = "  -3.81%  "
$ws.Range("E47").Value = "  -3.22%  "
$ws.Range("D48").Value = "1.43"
$ws.Range("E48").Value = "  -2.61%  "
$ws.Range("E49").Value = "  -1.06%  "
$ws.Range("D51").Value = "7.33"
$ws.Range("E51").Value = "  -2.68%  "
